$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" year column (O) that mirrors the existing "2020"
# column (N): same row-by-row formatting plus the updated figures.
$ws.Range("N3:N14").Copy() | Out-Null
$ws.Range("O3:O14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 4 header: new year column
$ws.Range("O4").Value = 2021

# Data values for the new 2021 column
$ws.Range("O5").Value = 97
$ws.Range("O6").Value = 96.2
$ws.Range("O7").Value = 62.7
$ws.Range("O8").Value = 100
$ws.Range("O9").Value = 100
$ws.Range("O10").Value = "-"
$ws.Range("O11").Value = 100
$ws.Range("O12").Value = 57.9
$ws.Range("O13").Value = 100
$ws.Range("O14").Value = "-"

# O3 stays blank (just formatted like N3 from the style copy above)

# Update the selected cell to mirror the recorded selection in the sheet view
$ws.Range("O17").Select() | Out-Null
